# Auto-generated: updates Price (D) and Volume(1h) (E) columns for the cryptos
# worksheet to reflect the latest scrape (GitHub Actions "Updated cryptos list" run).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.931.23"
$ws.Range("E2").Value = "  -2.05%  "

$ws.Range("D3").Value = "1.901.77"
$ws.Range("E3").Value = "  -4.09%  "

$ws.Range("D4").Value = "'1.005"
$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("E5").Value = "  -1.07%  "

$ws.Range("E6").Value = "  -0.08%  "

$ws.Range("E7").Value = "  -1.44%  "

$ws.Range("D8").Value = "'0.3811"
$ws.Range("E8").Value = "  -2.78%  "

$ws.Range("D9").Value = "'0.07702"
$ws.Range("E9").Value = "  -3.11%  "

$ws.Range("D10").Value = "'0.9736"
$ws.Range("E10").Value = "  -2.06%  "

$ws.Range("E11").Value = "  -4.00%  "

$ws.Range("D12").Value = "1.933.16"
$ws.Range("E12").Value = "  -2.36%  "

$ws.Range("D13").Value = "'6.919"
$ws.Range("E13").Value = "  -3.84%  "

$ws.Range("D14").Value = "'5.646"
$ws.Range("E14").Value = "  -3.17%  "

$ws.Range("D15").Value = "'0.07069"
$ws.Range("E15").Value = "  -0.43%  "

$ws.Range("E16").Value = "  -0.13%  "

$ws.Range("D17").Value = "'83.88"
$ws.Range("E17").Value = "  -4.35%  "

$ws.Range("D18").Value = "'0.000009475"
$ws.Range("E18").Value = "  -4.93%  "

$ws.Range("D19").Value = "'16.60"
$ws.Range("E19").Value = "  -4.16%  "

$ws.Range("E20").Value = "  -0.02%  "

$ws.Range("D21").Value = "28.902.97"
$ws.Range("E21").Value = "  -2.18%  "

$ws.Range("D22").Value = "'5.275"
$ws.Range("E22").Value = "  -5.33%  "

$ws.Range("D23").Value = "'10.85"
$ws.Range("E23").Value = "  -2.99%  "

$ws.Range("E24").Value = "  -0.60%  "

$ws.Range("D25").Value = "'157.94"
$ws.Range("E25").Value = "  -0.54%  "

$ws.Range("D26").Value = "'18.99"
$ws.Range("E26").Value = "  -3.17%  "

$ws.Range("D27").Value = "'5.619"
$ws.Range("E27").Value = "  -4.04%  "

$ws.Range("D28").Value = "'117.41"
$ws.Range("E28").Value = "  -1.84%  "

$ws.Range("D29").Value = "'1.836"
$ws.Range("E29").Value = "  -3.40%  "

$ws.Range("D30").Value = "'0.09251"
$ws.Range("E30").Value = "  -1.88%  "

$ws.Range("D31").Value = "'0.8567"
$ws.Range("E31").Value = "  -4.22%  "

$ws.Range("D32").Value = "'5.074"
$ws.Range("E32").Value = "  -3.13%  "

$ws.Range("D33").Value = "'1.237"
$ws.Range("E33").Value = "  -6.77%  "

$ws.Range("D34").Value = "'2.940"
$ws.Range("E34").Value = "  -8.06%  "

$ws.Range("D35").Value = "'0.05673"
$ws.Range("E35").Value = "  -2.39%  "

$ws.Range("D36").Value = "'1.139"
$ws.Range("E36").Value = "  -3.21%  "

$ws.Range("E37").Value = "  -0.01%  "

$ws.Range("D38").Value = "'0.02030"
$ws.Range("E38").Value = "  -3.22%  "

$ws.Range("D39").Value = "'0.5477"
$ws.Range("E39").Value = "  -4.60%  "

$ws.Range("D40").Value = "'7.379"
$ws.Range("E40").Value = "  -5.91%  "

$ws.Range("D41").Value = "'0.1752"
$ws.Range("E41").Value = "  -3.02%  "

$ws.Range("D42").Value = "'9.272"
$ws.Range("E42").Value = "  -4.34%  "

$ws.Range("D43").Value = "'2.755"
$ws.Range("E43").Value = "  -1.47%  "

$ws.Range("D44").Value = "'0.5153"
$ws.Range("E44").Value = "  -4.00%  "

$ws.Range("D45").Value = "'11.21"
$ws.Range("E45").Value = "  -5.55%  "

$ws.Range("D46").Value = "'0.06820"
$ws.Range("E46").Value = "  -1.80%  "

$ws.Range("D47").Value = "'2.057"
$ws.Range("E47").Value = "  -5.32%  "

$ws.Range("D48").Value = "'0.000002579"
$ws.Range("E48").Value = "  -16.54%  "

$ws.Range("D49").Value = "'110.13"
$ws.Range("E49").Value = "  -3.53%  "

$ws.Range("D50").Value = "'1.765"
$ws.Range("E50").Value = "  -3.45%  "

$ws.Range("D51").Value = "'1.004"
$ws.Range("E51").Value = "  -0.06%  "
